# Actualizacion automatica del mapa (2025-10-31 15:06:15)
# One incident record (Caso -28, "ZABALA /ALT/ 2836") was removed from the
# INCO sheet. Remove spreadsheet row 5 entirely so every following row
# shifts up by one and the sheet's used range shrinks from R35 to R34.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("INCO")

$ws.Rows.Item(5).EntireRow.Delete()
